$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 409 ("「彼を覚えていますか」...") entirely; rows below shift up by one.
$ws.Rows.Item(409).Delete()
